$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: "Subsequent biweekly meetings: 11:00 am on Monday." becomes
#           four separate runs: "Subsequent meetings: 11:00 am on Monday"
#           + "s " + "biweekly" + "."
# ---------------------------------------------------------------------

# Locate the target sentence robustly via Find (works regardless of any
# preceding edits / position shifts).
$findRng = $d.Content
$gotIt = $findRng.Find.Execute("Subsequent biweekly meetings: 11:00 am on Monday.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($gotIt) {
    $s0 = $findRng.Start

    # Remove "biweekly " (11 chars of "Subsequent " precede it, 9 chars long).
    $delRange = $d.Range($s0 + 11, $s0 + 20)
    $delRange.Text = ""

    # Text is now "Subsequent meetings: 11:00 am on Monday." starting at $s0.
    $newText = "Subsequent meetings: 11:00 am on Monday."
    $mondayEnd = $s0 + $newText.IndexOf("Monday") + ("Monday").Length

    # Insert "s " right after "Monday".
    $insS = $d.Range($mondayEnd, $mondayEnd)
    $insS.InsertAfter("s ")

    # Insert "biweekly" right before the trailing period.
    $afterMondayS = "Subsequent meetings: 11:00 am on Mondays ."
    $periodPos = $s0 + $afterMondayS.Length - 1
    $insBiweekly = $d.Range($periodPos, $periodPos)
    $insBiweekly.InsertAfter("biweekly")

    # At this point the engine has coalesced all the inserted text back
    # into a single run (adjacent runs with identical formatting are
    # merged automatically on every edit). Force the four required run
    # boundaries by dropping (and immediately removing) bookmarks at
    # each split point -- a bookmark marker splits the underlying run
    # even though it leaves no trace once removed.
    $b1 = $mondayEnd                    # Monday | s
    $b2 = $mondayEnd + 2                # s  | biweekly
    $b3 = $b2 + ("biweekly").Length     # biweekly | .

    $bm1 = "zzSplitBm1"
    $bm2 = "zzSplitBm2"
    $bm3 = "zzSplitBm3"
    $d.Bookmarks.Add($bm1, $d.Range($b1, $b1))
    $d.Bookmarks.Add($bm2, $d.Range($b2, $b2))
    $d.Bookmarks.Add($bm3, $d.Range($b3, $b3))
    $d.Bookmarks($bm1).Delete()
    $d.Bookmarks($bm2).Delete()
    $d.Bookmarks($bm3).Delete()
}

# ---------------------------------------------------------------------
# Change 2: the "Date:" paragraph's two adjacent runs "February" and
#           " 2026" (identical rPr: Arial / color 222222 / white shade)
#           merge into a single run "February 2026".
# ---------------------------------------------------------------------
$dateRange = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Date:*11:00 am*Tuesday*February*") {
        $dateRange = $d.Range($p.Range.Start, $p.Range.End)
    }
}

if ($dateRange -ne $null) {
    $null = $dateRange.Find.Execute("February 2026", $true, $false, $false, $false, $false, $true, 1, $false, "February 2026", 2)
}

Write-Output "done"
